$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values in columns V (TOTAL_DEATHS_BY_DOD) and W (NEW_DEATHS_BY_DOD)
# --- for rows 295-364, reflecting the switch of data source to CDC. ---
$ws.Range("V295").Value = 7033
$ws.Range("W295").Value = 107
$ws.Range("V296").Value = 7131
$ws.Range("V297").Value = 7236
$ws.Range("V298").Value = 7316
$ws.Range("V299").Value = 7408
$ws.Range("V300").Value = 7496
$ws.Range("V301").Value = 7586
$ws.Range("V302").Value = 7684
$ws.Range("V303").Value = 7783
$ws.Range("V304").Value = 7882
$ws.Range("W304").Value = 99
$ws.Range("V305").Value = 7980
$ws.Range("V306").Value = 8069
$ws.Range("V307").Value = 8173
$ws.Range("V308").Value = 8281
$ws.Range("V309").Value = 8396
$ws.Range("V310").Value = 8512
$ws.Range("V311").Value = 8602
$ws.Range("V312").Value = 8718
$ws.Range("V313").Value = 8807
$ws.Range("V314").Value = 8914
$ws.Range("V315").Value = 9016
$ws.Range("V316").Value = 9114
$ws.Range("V317").Value = 9191
$ws.Range("V318").Value = 9297
$ws.Range("V319").Value = 9389
$ws.Range("V320").Value = 9476
$ws.Range("V321").Value = 9566
$ws.Range("V322").Value = 9647
$ws.Range("V323").Value = 9746
$ws.Range("V324").Value = 9832
$ws.Range("V325").Value = 9940
$ws.Range("V326").Value = 10015
$ws.Range("V327").Value = 10086
$ws.Range("V328").Value = 10169
$ws.Range("V329").Value = 10244
$ws.Range("V330").Value = 10327
$ws.Range("V331").Value = 10389
$ws.Range("V332").Value = 10460
$ws.Range("V333").Value = 10515
$ws.Range("V334").Value = 10565
$ws.Range("V335").Value = 10613
$ws.Range("V336").Value = 10664
$ws.Range("V337").Value = 10707
$ws.Range("V338").Value = 10746
$ws.Range("V339").Value = 10795
$ws.Range("V340").Value = 10842
$ws.Range("V341").Value = 10900
$ws.Range("V342").Value = 10926
$ws.Range("V343").Value = 10962
$ws.Range("V344").Value = 11001
$ws.Range("V345").Value = 11035
$ws.Range("W345").Value = 34
$ws.Range("V346").Value = 11069
$ws.Range("V347").Value = 11095
$ws.Range("V348").Value = 11127
$ws.Range("W349").Value = 28
$ws.Range("V351").Value = 11211
$ws.Range("W351").Value = 30
$ws.Range("V352").Value = 11248
$ws.Range("V353").Value = 11277
$ws.Range("W353").Value = 29
$ws.Range("V354").Value = 11305
$ws.Range("V355").Value = 11323
$ws.Range("V356").Value = 11347
$ws.Range("V357").Value = 11362
$ws.Range("V358").Value = 11378
$ws.Range("V359").Value = 11388
$ws.Range("W359").Value = 10
$ws.Range("V360").Value = 11401
$ws.Range("W360").Value = 13
$ws.Range("V361").Value = 11409
$ws.Range("V362").Value = 11416
$ws.Range("W362").Value = 7
$ws.Range("V363").Value = 11419
$ws.Range("W363").Value = 3
$ws.Range("V364").Value = 11424
$ws.Range("W364").Value = 5

# --- Append new row 365 (2021-03-02, serial date 44257) ---
# Copy date formatting from A364 into A365 before setting its value so the new
# row reuses the existing date style (rather than minting a new numFmt/style).
$ws.Range("A364").Copy()
$ws.Range("A365").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A365").Value = 44257
$ws.Range("B365").Value = 776337
$ws.Range("C365").Value = 644
$ws.Range("D365").Value = 650084
$ws.Range("E365").Value = 420
$ws.Range("F365").Value = 126253
$ws.Range("G365").Value = 224
$ws.Range("H365").Value = 749637
$ws.Range("I365").Value = 6047838
$ws.Range("J365").Value = 6797475
$ws.Range("K365").Value = 7505
$ws.Range("L365").Value = 15
$ws.Range("M365").Value = 11436
$ws.Range("P365").Value = -392
$ws.Range("Q365").Value = 13125
$ws.Range("R365").Value = 1021
$ws.Range("S365").Value = 751776
$ws.Range("T365").Value = 59
$ws.Range("U365").Value = 18679
$ws.Range("V365").Value = 11426
$ws.Range("W365").Value = 2

# --- Update the workbook-level defined name / print range to include row 365 ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "ALL_DAILY_CASE_INFO_PUBLIC") {
        $n.RefersTo = "=ALL_DAILY_CASE_INFO_PUBLIC!`$A`$1:`$W`$365"
    }
}
